# Collection-data scripts adapted for parallelization: append page 16
# (4 new specimen entries) to the HJ-7 raw-data collection entry sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared values repeated across the whole new batch (page 16, 1981-07-16,
# Galiano Island, Canada / British Columbia, curation remark + locality
# remark already used by the prior rows on this sheet).
$conf       = "h"
$theDate    = 19810716
$locality   = "Galiano Island"
$country    = "Canada"
$state      = "British Columbia"
$island     = "Galiano Island"
$curation   = "red dot crossed out next to taxon name"
$locRemarks = "set of collections from Montague to Georgeson Bay and up to Sphagnum bog of Gauner Road"

$newRows = @(
    @{ Row=116; NumPage=1; RecordNum=2043; VName="Matthiola incana";           VSciName="Matthiola incana";           SciName="Matthiola incana" },
    @{ Row=117; NumPage=2; RecordNum=2044; VName="Polygonum aviculare";        VSciName="Polygonum aviculare";        SciName="Polygonum aviculare" },
    @{ Row=118; NumPage=3; RecordNum=2045; VName="Actea rubra";                VSciName="Actaea rubra";                SciName="Actaea rubra" },
    @{ Row=119; NumPage=4; RecordNum=2046; VName="Polygonum spergulariiforme"; VSciName="Polygonum spergulariiforme"; SciName="Polygonum spergulariiforme" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = 16             # A pageNum
    $ws.Cells.Item($row, 2).Value  = $r.NumPage     # B numPage
    $ws.Cells.Item($row, 3).Value  = $r.RecordNum   # C recordNum
    $ws.Cells.Item($row, 4).Value  = $r.VName       # D vName
    $ws.Cells.Item($row, 5).Value  = $r.VSciName    # E vSciName
    $ws.Cells.Item($row, 6).Value  = $conf          # F conf
    $ws.Cells.Item($row, 7).Value  = $r.SciName     # G sciName
    $ws.Cells.Item($row, 8).Value  = $theDate       # H date
    $ws.Cells.Item($row, 9).Value  = $locality      # I locality
    $ws.Cells.Item($row, 10).Value = $country       # J country
    $ws.Cells.Item($row, 11).Value = $state         # K stateProvince
    $ws.Cells.Item($row, 12).Value = $island        # L island
    $ws.Cells.Item($row, 13).Value = $curation      # M curationMetadata
    $ws.Cells.Item($row, 17).Value = $locRemarks    # Q locationRemarks

    $ws.Rows.Item($row).RowHeight = 102
}

# Match the committed view state: zoomed to 134%, scrolled/frozen near the
# new rows, with the bottom-left pane's active cell on the new D117 entry.
$excel.ActiveWindow.Zoom = 134
$null = $ws.Range("A115").Select()
$null = $ws.Range("D117").Select()
